$wb = $excel.ActiveWorkbook
$wsRunManager = $wb.Worksheets.Item("RUN_MANAGER")
$wsTestData = $wb.Worksheets.Item("TEST_DATA")

# --- RUN_MANAGER sheet: flip Execute flag to "yes" for rows 2, 3, 5 ---
$wsRunManager.Range("D2").Value = "yes"
$wsRunManager.Range("D3").Value = "yes"
$wsRunManager.Range("D5").Value = "yes"

# --- TEST_DATA sheet: switch Execute/Browser for rows 2,3,5,6,20 to no/edge ---
$wsTestData.Range("B2").Value = "no"
$wsTestData.Range("C2").Value = "edge"

$wsTestData.Range("B3").Value = "no"
$wsTestData.Range("C3").Value = "edge"

$wsTestData.Range("B5").Value = "no"
$wsTestData.Range("C5").Value = "edge"

$wsTestData.Range("B6").Value = "no"
$wsTestData.Range("C6").Value = "edge"

$wsTestData.Range("B20").Value = "no"
$wsTestData.Range("C20").Value = "edge"

# --- Update selections / view state to match the saved worksheet view ---
$wsRunManager.Activate() | Out-Null
$wsRunManager.Range("D3").Select() | Out-Null

$wsTestData.Activate() | Out-Null
$wsTestData.Range("B3").Select() | Out-Null
